$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1315.4231
$ws.Range("J17").Value = 1300.0714
$ws.Range("L17").Value = 3900.2142
$ws.Range("N17").Value = -4236.2142
$ws.Range("H33").Value = 110.478264
$ws.Range("I33").Value = 50.055557
$ws.Range("J33").Value = 328
$ws.Range("K33").Value = 50.055557
$ws.Range("L33").Value = 328
$ws.Range("M33").Value = 178.944443
$ws.Range("N33").Value = -786
$ws.Range("H116").Value = 442135.47
$ws.Range("I116").Value = 1251938.8
$ws.Range("J116").Value = 10240.4
$ws.Range("K116").Value = 1251938.8
$ws.Range("L116").Value = 10240.4
$ws.Range("M116").Value = -1248496.8
$ws.Range("N116").Value = -17124.4
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H125").Value = 4066.8333
$ws.Range("I125").Value = 3244.4
$ws.Range("J125").Value = 4654.2856
$ws.Range("K125").Value = 29199.6
$ws.Range("L125").Value = 41888.5704
$ws.Range("M125").Value = -26739.6
$ws.Range("N125").Value = -46808.5704
$ws.Range("H132").Value = 6404.7715
$ws.Range("I132").Value = 6360.577
$ws.Range("J132").Value = 6532.4443
$ws.Range("K132").Value = 19081.731
$ws.Range("L132").Value = 19597.3329
$ws.Range("M132").Value = -16551.731
$ws.Range("N132").Value = -24657.3329
$ws.Range("H141").Value = 69305.734
$ws.Range("I141").Value = 93063.27
$ws.Range("K141").Value = 279189.81
$ws.Range("M141").Value = -274009.81
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8329.623
$ws.Range("I32").Value = 7131.75
$ws.Range("K32").Value = 7131.75
$ws.Range("M32").Value = -6844.75
$ws.Range("H110").Value = 1802.2
$ws.Range("I110").Value = 1003.6667
$ws.Range("K110").Value = 1003.6667
$ws.Range("M110").Value = 1041.3333
$ws.Range("H124").Value = 29429
$ws.Range("J124").Value = 29429
$ws.Range("L124").Value = 29429
$ws.Range("N124").Value = -39249
$ws.Range("H132").Value = 2268.8
$ws.Range("I132").Value = 1710.075
$ws.Range("J132").Value = 4503.7
$ws.Range("K132").Value = 5130.225
$ws.Range("L132").Value = 13511.1
$ws.Range("M132").Value = -2600.225
$ws.Range("N132").Value = -18571.1
$ws.Range("H137").Value = 39588.332
$ws.Range("J137").Value = 39588.332
$ws.Range("L137").Value = 39588.332
$ws.Range("N137").Value = -49788.332
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2422.1
$ws.Range("I105").Value = 2351.28
$ws.Range("J105").Value = 2776.2
$ws.Range("K105").Value = 2351.28
$ws.Range("L105").Value = 2776.2
$ws.Range("M105").Value = -604.2800000000002
$ws.Range("N105").Value = -6270.2
$ws.Range("H137").Value = 34920
$ws.Range("J137").Value = 39880
$ws.Range("L137").Value = 39880
$ws.Range("N137").Value = -50080
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13891473
$ws.Range("I16").Value = 27778944
$ws.Range("J16").Value = 4003.25
$ws.Range("K16").Value = 27778944
$ws.Range("L16").Value = 4003.25
$ws.Range("M16").Value = -27778657
$ws.Range("N16").Value = -4577.25
$ws.Range("H31").Value = 12198376
$ws.Range("I31").Value = 1867.8462
$ws.Range("K31").Value = 1867.8462
$ws.Range("M31").Value = -1572.8462
$ws.Range("H34").Value = 12198376
$ws.Range("I34").Value = 1867.8462
$ws.Range("K34").Value = 1867.8462
$ws.Range("M34").Value = -1665.8462
$ws.Range("H105").Value = 2159.625
$ws.Range("I105").Value = 1995.6666
$ws.Range("J105").Value = 2258
$ws.Range("K105").Value = 1995.6666
$ws.Range("L105").Value = 2258
$ws.Range("M105").Value = -248.6666
$ws.Range("N105").Value = -5752
$ws.Range("H113").Value = 13891473
$ws.Range("I113").Value = 27778944
$ws.Range("J113").Value = 4003.25
$ws.Range("K113").Value = 27778944
$ws.Range("L113").Value = 4003.25
$ws.Range("M113").Value = -27776774
$ws.Range("N113").Value = -8343.25
$ws.Range("H124").Value = 25225.2
$ws.Range("J124").Value = 25225.2
$ws.Range("L124").Value = 25225.2
$ws.Range("N124").Value = -30135.2
$ws.Range("H134").Value = 7535.905
$ws.Range("I134").Value = 15689.143
$ws.Range("K134").Value = 47067.429
$ws.Range("M134").Value = -44532.429
$ws.Range("H137").Value = 38920
$ws.Range("I137").Value = 9800
$ws.Range("J137").Value = 41346.668
$ws.Range("K137").Value = 9800
$ws.Range("L137").Value = 41346.668
$ws.Range("M137").Value = -4700
$ws.Range("N137").Value = -51546.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 788.10205
$ws.Range("I113").Value = 663.84
$ws.Range("J113").Value = 917.5417
$ws.Range("K113").Value = 1991.52
$ws.Range("L113").Value = 2752.6251
$ws.Range("M113").Value = 178.48
$ws.Range("N113").Value = -7092.6251
$ws.Range("H122").Value = 3646.86
$ws.Range("I122").Value = 613.25
$ws.Range("J122").Value = 3910.652
$ws.Range("K122").Value = 5519.25
$ws.Range("L122").Value = 35195.868
$ws.Range("M122").Value = -3069.25
$ws.Range("N122").Value = -40095.868
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5616.4385
$ws.Range("I70").Value = 5051.564
$ws.Range("J70").Value = 6840.3335
$ws.Range("K70").Value = 5051.564
$ws.Range("L70").Value = 6840.3335
$ws.Range("M70").Value = -4781.564
$ws.Range("N70").Value = -7380.3335
$ws.Range("H73").Value = 5616.4385
$ws.Range("I73").Value = 5051.564
$ws.Range("J73").Value = 6840.3335
$ws.Range("K73").Value = 5051.564
$ws.Range("L73").Value = 6840.3335
$ws.Range("M73").Value = -4115.564
$ws.Range("N73").Value = -8712.333500000001
$ws.Range("H113").Value = 1219
$ws.Range("I113").Value = 1017.75
$ws.Range("J113").Value = 1380
$ws.Range("K113").Value = 1017.75
$ws.Range("L113").Value = 1380
$ws.Range("M113").Value = 1152.25
$ws.Range("N113").Value = -5720
$ws.Range("H126").Value = 2955.04
$ws.Range("I126").Value = 2955.04
$ws.Range("K126").Value = 8865.119999999999
$ws.Range("M126").Value = -6395.119999999999
$ws.Range("H135").Value = 52540.5
$ws.Range("I135").Value = 139500
$ws.Range("J135").Value = 47425.234
$ws.Range("K135").Value = 139500
$ws.Range("L135").Value = 47425.234
$ws.Range("M135").Value = -134430
$ws.Range("N135").Value = -57565.234
$ws.Range("H136").Value = 24839.428
$ws.Range("J136").Value = 24839.428
$ws.Range("L136").Value = 74518.284
$ws.Range("N136").Value = -79618.284
$ws.Range("H137").Value = 42397
$ws.Range("J137").Value = 42397
$ws.Range("L137").Value = 42397
$ws.Range("N137").Value = -52597
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3509.742
$ws.Range("I7").Value = 1937.75
$ws.Range("J7").Value = 5186.533
$ws.Range("K7").Value = 1937.75
$ws.Range("L7").Value = 5186.533
$ws.Range("M7").Value = -1825.75
$ws.Range("N7").Value = -5410.533
$ws.Range("H16").Value = 1792.0667
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3340
$ws.Range("H61").Value = 1335.2667
$ws.Range("I61").Value = 1079.1538
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1079.1538
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -877.1538
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 1335.2667
$ws.Range("I113").Value = 1079.1538
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1079.1538
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1090.8462
$ws.Range("N113").Value = -7340
$ws.Range("H126").Value = 3509.742
$ws.Range("I126").Value = 1937.75
$ws.Range("J126").Value = 5186.533
$ws.Range("K126").Value = 5813.25
$ws.Range("L126").Value = 15559.599
$ws.Range("M126").Value = -3343.25
$ws.Range("N126").Value = -20499.599
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15158903
$ws.Range("I132").Value = 10909.728
$ws.Range("J132").Value = 30306896
$ws.Range("K132").Value = 32729.184
$ws.Range("L132").Value = 90920688
$ws.Range("M132").Value = -30199.184
$ws.Range("N132").Value = -90925748
